# Fix Issue #32: "Keycaps Filament" (row 9) and "Buttons" (row 10) had no
# "Stueck (in packung)" quantity entered yet, so the per-piece-cost formula
# (Kosten pro Stueck = D/B) divided by an empty/zero quantity and surfaced
# #DIV/0! errors. Filling in the missing quantities (and, for the Buttons
# row, the now-known per-package price) lets the shared formulas resolve.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: "Keycaps Filament" -------------------------------------------
# Quantity-per-package was never entered (division by blank -> #DIV/0!).
# One package, price still unknown/TBD (stays 0).
$ws.Range("B9").Value = 1
$ws.Range("D9").Value = 0

# --- Row 10: "Buttons" ----------------------------------------------------
# Same missing-quantity bug, plus this is where the actual button price
# (per piece, from the RS Components link in G10) now gets recorded.
$ws.Range("B10").Value = 1
$ws.Range("D10").Value = 0.138

# Recalculate so the shared E/F formulas (Kosten pro Stueck / Kosten Gesamt)
# and the downstream per-person total in B18 pick up the new inputs.
$excel.CalculateFull()

# --- Cosmetic follow-up from the same editing session --------------------
# Pre-format a block of rows further down the sheet (currency style on
# D:F, hyperlink style on G) the same way the existing product rows are
# styled, ready for more products to be appended later.
$ws.Range("D9:F9").Copy()
$ws.Range("D46:F48").PasteSpecial(-4122)
$ws.Range("G10").Copy()
$ws.Range("G46:G48").PasteSpecial(-4122)

$ws.Range("D9:F9").Copy()
$ws.Range("D49:F51").PasteSpecial(-4122)
$ws.Range("G10").Copy()
$ws.Range("G49:G51").PasteSpecial(-4122)

$ws.Range("D9:F9").Copy()
$ws.Range("D52:F52").PasteSpecial(-4122)

$ws.Range("D9:F9").Copy()
$ws.Range("D53:F53").PasteSpecial(-4122)
$ws.Range("G10").Copy()
$ws.Range("G53").PasteSpecial(-4122)

$ws.Range("D9:F9").Copy()
$ws.Range("D54:F56").PasteSpecial(-4122)

$ws.Range("D9").Copy()
$ws.Range("D57").PasteSpecial(-4122)

$ws.Range("D9").Copy()
$ws.Range("B61").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column widths were re-tightened (no more bestFit auto-sizing on most
# columns) after the extra formatted rows were added below the table.
$ws.Columns.Item(1).ColumnWidth = 61.28515625
$ws.Columns.Item(2).ColumnWidth = 13.5703125
$ws.Columns.Item(3).ColumnWidth = 9.5703125
$ws.Columns.Item(4).ColumnWidth = 13.28515625
$ws.Columns.Item(5).ColumnWidth = 14.85546875
$ws.Columns.Item(6).ColumnWidth = 9.5703125
$ws.Columns.Item(7).ColumnWidth = 56.85546875
$ws.Columns.Item(8).ColumnWidth = 23.7109375

# View state: zoomed in a bit more and left on a different active cell.
$excel.ActiveWindow.Zoom = 145
$ws.Range("E15").Select()

$excel.CalculateFull()
